# Fix "appartments" typo -> "apartments" on Sheet3 and Sheet4, tidy up the
# redundant cell format left behind on Sheet4!A4, re-point the saved
# selection/active-sheet state (Sheet2 becomes the active tab/cell), and
# nudge Sheet4's first column width by a rounding tick.

$wb = $excel.ActiveWorkbook

# --- Sheet1: move the saved selection from G14 (out of range) to G10 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G10").Select()

# --- Sheet2: becomes the active sheet/tab, selection moves to B2 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("B2").Select()

# --- Sheet3: fix "appartments" typo, move saved selection to G4 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Total of apartments"
$ws3.Range("G4").Select()

# --- Sheet4: fix "appartments " typo (trailing space kept) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "Total of apartments "

# Sheet4!A4 had a one-off duplicate font (Calibri 12) that is identical to
# the workbook's base font; reset it back to the shared/default font so the
# redundant font entry is no longer referenced.
$ws4.Range("A4").Font.Name = "Calibri"
$ws4.Range("A4").Font.Size = 12

# Column A on Sheet4 widens very slightly (17.85 -> 17.86 characters)
$ws4.Columns.Item(1).ColumnWidth = 17.86

# Re-select B4 (kept the same) and make sure Sheet2 ends up the active tab
$ws4.Range("B4").Select()
$ws2.Activate()
$ws2.Range("B2").Select()

Write-Output "done"
